$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.14230000000001
$ws.Range("A10").Value = -20.47339999999997
$ws.Range("A12").Value = -22.41370000000003
$ws.Range("E13").Value = 12.2045
$ws.Range("A18").Value = -22.20900000000002
